$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new log entry row (row 12) with the activity text and hours worked.
$ws.Range("B12").Value = "Tackled the methodologies section"
$ws.Range("C12").Value = 0.75

# Make sure the SUM formula in F1 still covers the new row and recalculates.
$ws.Range("F1").Formula = "=SUM(C2:C32)"
$excel.Calculate()

# Update the view state to match what was saved: scrolled down a bit with F10 selected.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("F10").Select()

$wb.Save()
